$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 73, shifting existing rows 73:141 down to 74:142.
$ws.Rows("73:73").Insert()

# Populate the newly inserted row 73 with the new record (mirrors the
# constant columns shared by every other data row, plus the new values).
$ws.Range("A73").Value = 8
$ws.Range("B73").Value = "Terminal La Palmera de La Serena"
$ws.Range("C73").Value = "Coquimbo"
$ws.Range("D73").Value = 44729
$ws.Range("E73").Value = 4
$ws.Range("F73").Value = 100112001
$ws.Range("G73").Value = "Berenjena"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 460
$ws.Range("K73").Value = 8000
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = 8500
$ws.Range("N73").Value = "$/caja 50 unidades"
$ws.Range("O73").Value = "Región de Arica y Parinacota"
$ws.Range("P73").Value = 170
$ws.Range("Q73").Value = 50
$ws.Range("R73").Value = "Hortaliza"
